# Sync attendance_reports: swap "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# in column G ("Recorded By") wherever it matches exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $currentValue = $cell.Value()
    if ($currentValue -eq $oldValue) {
        $cell.Value = $newValue
    }
}
